# Applies the Fri Aug 11 22:24:52 UTC 2023 cryptos-list refresh:
#  - updates Price (D) / Volume 1h (E) figures for existing rows
#  - drops BabyDogeCoin (old row 46); Aptos..Algorand shift up one row
#  - adds a new Cronos row at the bottom (row 51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.394.24"
$ws.Cells.Item(2, 5).Value = "  -0.10%  "

$ws.Cells.Item(3, 4).Value = "1.844.27"
$ws.Cells.Item(3, 5).Value = "  -0.27%  "

$ws.Cells.Item(5, 4).Value = "'238.87"
$ws.Cells.Item(5, 5).Value = "  -0.87%  "

$ws.Cells.Item(6, 4).Value = "'0.6311"
$ws.Cells.Item(6, 5).Value = "  -0.53%  "

$ws.Cells.Item(8, 4).Value = "'0.07532"
$ws.Cells.Item(8, 5).Value = "  -0.41%  "

$ws.Cells.Item(9, 4).Value = "'0.2930"
$ws.Cells.Item(9, 5).Value = "  -1.38%  "

$ws.Cells.Item(10, 4).Value = "'24.53"
$ws.Cells.Item(10, 5).Value = "  -0.50%  "

$ws.Cells.Item(11, 4).Value = "'0.07712"
$ws.Cells.Item(11, 5).Value = "  -0.08%  "

$ws.Cells.Item(12, 4).Value = "1.830.36"
$ws.Cells.Item(12, 5).Value = "  -7.79%  "

$ws.Cells.Item(13, 4).Value = "'5.001"
$ws.Cells.Item(13, 5).Value = "  +0.12%  "

$ws.Cells.Item(14, 4).Value = "'0.6798"
$ws.Cells.Item(14, 5).Value = "  -0.97%  "

$ws.Cells.Item(15, 4).Value = "'0.00001044"
$ws.Cells.Item(15, 5).Value = "  +5.17%  "

$ws.Cells.Item(16, 4).Value = "'83.35"
$ws.Cells.Item(16, 5).Value = "  +0.38%  "

$ws.Cells.Item(17, 4).Value = "2.092.52"
$ws.Cells.Item(17, 5).Value = "  -7.59%  "

$ws.Cells.Item(18, 4).Value = "'6.171"
$ws.Cells.Item(18, 5).Value = "  -0.70%  "

$ws.Cells.Item(19, 4).Value = "29.424.97"
$ws.Cells.Item(19, 5).Value = "  -0.12%  "

$ws.Cells.Item(20, 4).Value = "'228.55"
$ws.Cells.Item(20, 5).Value = "  -1.55%  "

$ws.Cells.Item(21, 5).Value = "  -0.71%  "

$ws.Cells.Item(22, 5).Value = "  +0.03%  "

$ws.Cells.Item(23, 5).Value = "  -1.94%  "

$ws.Cells.Item(24, 5).Value = "  +0.04%  "

$ws.Cells.Item(25, 4).Value = "'156.99"
$ws.Cells.Item(25, 5).Value = "  +0.72%  "

$ws.Cells.Item(26, 4).Value = "'0.1393"
$ws.Cells.Item(26, 5).Value = "  +0.39%  "

$ws.Cells.Item(27, 4).Value = "'8.354"
$ws.Cells.Item(27, 5).Value = "  -0.80%  "

$ws.Cells.Item(28, 5).Value = "  -0.58%  "

$ws.Cells.Item(29, 4).Value = "'1.458"
$ws.Cells.Item(29, 5).Value = "  -0.56%  "

$ws.Cells.Item(30, 4).Value = "'1.282"
$ws.Cells.Item(30, 5).Value = "  +1.76%  "

$ws.Cells.Item(31, 4).Value = "'0.05630"
$ws.Cells.Item(31, 5).Value = "  -3.14%  "

$ws.Cells.Item(32, 5).Value = "  -0.87%  "

$ws.Cells.Item(33, 4).Value = "'4.028"
$ws.Cells.Item(33, 5).Value = "  +0.13%  "

$ws.Cells.Item(34, 4).Value = "'1.846"
$ws.Cells.Item(34, 5).Value = "  -0.78%  "

$ws.Cells.Item(35, 4).Value = "'1.156"
$ws.Cells.Item(35, 5).Value = "  -0.14%  "

$ws.Cells.Item(36, 4).Value = "'0.7125"
$ws.Cells.Item(36, 5).Value = "  -0.66%  "

$ws.Cells.Item(37, 4).Value = "'2.590"
$ws.Cells.Item(37, 5).Value = "  +0.07%  "

$ws.Cells.Item(38, 4).Value = "1.247.17"
$ws.Cells.Item(38, 5).Value = "  -0.41%  "

$ws.Cells.Item(39, 5).Value = "  +0.14%  "

$ws.Cells.Item(40, 4).Value = "'2.765"
$ws.Cells.Item(40, 5).Value = "  -1.21%  "

$ws.Cells.Item(41, 4).Value = "'6.375"
$ws.Cells.Item(41, 5).Value = "  +4.63%  "

$ws.Cells.Item(42, 4).Value = "'0.9027"
$ws.Cells.Item(42, 5).Value = "  +0.04%  "

$ws.Cells.Item(44, 4).Value = "'101.65"
$ws.Cells.Item(44, 5).Value = "  -0.12%  "

$ws.Cells.Item(45, 4).Value = "'65.76"
$ws.Cells.Item(45, 5).Value = "  -1.97%  "

# Row 46: "BabyDogeCoin" -> "Aptos"
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).Value = "'7.098"
$ws.Cells.Item(46, 5).Value = "  -1.62%  "

# Row 47: "Aptos" -> "TheSandbox"
$ws.Cells.Item(47, 2).Value = "TheSandbox"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(47, 4).Value = "'0.3994"
$ws.Cells.Item(47, 5).Value = "  -0.64%  "

# Row 48: "TheSandbox" -> "EnergySwap"
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).Value = "'8.924"
$ws.Cells.Item(48, 5).Value = "  -2.43%  "

# Row 49: "EnergySwap" -> "RenderToken"
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).Value = "'1.672"
$ws.Cells.Item(49, 5).Value = "  -0.83%  "

# Row 50: "RenderToken" -> "Algorand"
$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).Value = "'0.1122"
$ws.Cells.Item(50, 5).Value = "  -0.41%  "

# Row 51: "Algorand" -> "Cronos"
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.05705"
$ws.Cells.Item(51, 5).Value = "  -0.80%  "
